$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C3").Value = "[0.016612939351596613, 0.07334226991833549, 0.07541311194742319]"
$ws.Range("E3").Value = "[5.008104595025083, 5.734199448081883, 10.59963755938173]"

$ws.Range("C5").Value = "[0.016612939351596613, 0.07334226991833549, 0.07541311194742319]"
$ws.Range("E5").Value = "[5.00810459500194, 5.734199448076016, 638.903418288658]"

$ws.Range("C7").Value = "[0.016612939351596613, 0.07334226991833549, 0.07541311194742319]"
$ws.Range("E7").Value = "[5.008104595020174, 5.734199448068993, 10.599637559494415]"
